$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.672209024429321
$ws.Range("B1").Value = 1.594421029090881
$ws.Range("C1").Value = 4.751124382019043
$ws.Range("D1").Value = 1.208233952522278
$ws.Range("E1").Value = 0.6363009214401245
